# Fireball Test workbook update
# - Activate sheet "S22" (was "S66")
# - Update selection on S22 to G3
# - Resize several columns on S22, add width for column G
# - Add "Fireball test2" label in G1
# - Fill "x" markers in K4, K5, K6, K8

$wb = $excel.ActiveWorkbook

$ws22 = $wb.Worksheets.Item("S22")

# Make S22 the active sheet/tab (S66 was previously active)
$ws22.Activate()

# Column width changes on S22 (values chosen so the engine's internal
# rounding lands as close as possible on the authored widths)
$ws22.Columns.Item(2).ColumnWidth = 27
$ws22.Columns.Item(4).ColumnWidth = 19.666666666666668
$ws22.Columns.Item(5).ColumnWidth = 14
$ws22.Columns.Item(6).ColumnWidth = 14.666666666666666
$ws22.Columns.Item(7).ColumnWidth = 15.166666666666666
$ws22.Columns.Item(8).ColumnWidth = 12.833333333333334

# New header text and "x" markers
$ws22.Range("G1").Value = "Fireball test2"
$ws22.Range("K4").Value = "x"
$ws22.Range("K5").Value = "x"
$ws22.Range("K6").Value = "x"
$ws22.Range("K8").Value = "x"

# Selection on S22 moves from K4 to G3
$ws22.Range("G3").Select()
